$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 46073
$ws.Cells.Item(2, 2).Value = 9.81
$ws.Cells.Item(2, 3).Value = 47.478
$ws.Cells.Item(2, 4).Value = 47.478
$ws.Cells.Item(2, 6).Value = "20.02.20261"
$ws.Cells.Item(3, 1).Value = 46073.01041666666
$ws.Cells.Item(3, 2).Value = 7.581
$ws.Cells.Item(3, 3).Value = 117.145
$ws.Cells.Item(3, 4).Value = 117.145
$ws.Cells.Item(3, 6).Value = "20.02.20262"
$ws.Cells.Item(4, 1).Value = 46073.02083333334
$ws.Cells.Item(4, 2).Value = -36.616
$ws.Cells.Item(4, 3).Value = 602.662
$ws.Cells.Item(4, 4).Value = 602.662
$ws.Cells.Item(4, 6).Value = "20.02.20263"
$ws.Cells.Item(5, 1).Value = 46073.03125
$ws.Cells.Item(5, 2).Value = -42.395
$ws.Cells.Item(5, 3).Value = 608.649
$ws.Cells.Item(5, 4).Value = 608.649
$ws.Cells.Item(5, 6).Value = "20.02.20264"
$ws.Cells.Item(6, 1).Value = 46073.04166666666
$ws.Cells.Item(6, 2).Value = -26.475
$ws.Cells.Item(6, 3).Value = 778.893
$ws.Cells.Item(6, 4).Value = 778.893
$ws.Cells.Item(6, 6).Value = "20.02.20265"
$ws.Cells.Item(7, 1).Value = 46073.05208333334
$ws.Cells.Item(7, 2).Value = -15.193
$ws.Cells.Item(7, 3).Value = 655.303
$ws.Cells.Item(7, 4).Value = 655.303
$ws.Cells.Item(7, 6).Value = "20.02.20266"
$ws.Cells.Item(8, 1).Value = 46073.0625
$ws.Cells.Item(8, 2).Value = -46.129
$ws.Cells.Item(8, 3).Value = 667.134
$ws.Cells.Item(8, 4).Value = 667.134
$ws.Cells.Item(8, 6).Value = "20.02.20267"
$ws.Cells.Item(9, 1).Value = 46073.07291666666
$ws.Cells.Item(9, 2).Value = -18.88
$ws.Cells.Item(9, 3).Value = 725.826
$ws.Cells.Item(9, 4).Value = 725.826
$ws.Cells.Item(9, 6).Value = "20.02.20268"
$ws.Cells.Item(10, 1).Value = 46073.08333333334
$ws.Cells.Item(10, 2).Value = 14.727
$ws.Cells.Item(10, 3).Value = 102.864
$ws.Cells.Item(10, 4).Value = 102.864
$ws.Cells.Item(10, 6).Value = "20.02.20269"
$ws.Cells.Item(11, 1).Value = 46073.09375
$ws.Cells.Item(11, 2).Value = 36.674
$ws.Cells.Item(11, 3).Value = -363.096
$ws.Cells.Item(11, 4).Value = -363.096
$ws.Cells.Item(11, 6).Value = "20.02.202610"
$ws.Cells.Item(12, 1).Value = 46073.10416666666
$ws.Cells.Item(12, 2).Value = 30.492
$ws.Cells.Item(12, 3).Value = -158.658
$ws.Cells.Item(12, 4).Value = -158.658
$ws.Cells.Item(12, 6).Value = "20.02.202611"
$ws.Cells.Item(13, 1).Value = 46073.11458333334
$ws.Cells.Item(13, 2).Value = 16.571
$ws.Cells.Item(13, 3).Value = 76.65900000000001
$ws.Cells.Item(13, 4).Value = 76.65900000000001
$ws.Cells.Item(13, 6).Value = "20.02.202612"
$ws.Cells.Item(14, 1).Value = 46073.125
$ws.Cells.Item(14, 2).Value = 32.099
$ws.Cells.Item(14, 3).Value = -13.381
$ws.Cells.Item(14, 4).Value = -13.381
$ws.Cells.Item(14, 6).Value = "20.02.202613"
$ws.Cells.Item(15, 1).Value = 46073.13541666666
$ws.Cells.Item(15, 2).Value = 66.279
$ws.Cells.Item(15, 3).Value = -254.32
$ws.Cells.Item(15, 4).Value = -254.32
$ws.Cells.Item(15, 6).Value = "20.02.202614"
$ws.Cells.Item(16, 1).Value = 46073.14583333334
$ws.Cells.Item(16, 2).Value = 9.887
$ws.Cells.Item(16, 3).Value = 19.667
$ws.Cells.Item(16, 4).Value = 19.667
$ws.Cells.Item(16, 6).Value = "20.02.202615"
$ws.Cells.Item(17, 1).Value = 46073.15625
$ws.Cells.Item(17, 2).Value = -25.201
$ws.Cells.Item(17, 3).Value = 550
$ws.Cells.Item(17, 4).Value = 550
$ws.Cells.Item(17, 6).Value = "20.02.202616"
$ws.Cells.Item(18, 1).Value = 46073.16666666666
$ws.Cells.Item(18, 2).Value = -17.062
$ws.Cells.Item(18, 3).Value = 611.774
$ws.Cells.Item(18, 4).Value = 611.774
$ws.Cells.Item(18, 6).Value = "20.02.202617"
$ws.Cells.Item(19, 1).Value = 46073.17708333334
$ws.Cells.Item(19, 2).Value = -0.804
$ws.Cells.Item(19, 3).Value = 611.7910000000001
$ws.Cells.Item(19, 4).Value = 611.7910000000001
$ws.Cells.Item(19, 6).Value = "20.02.202618"
$ws.Cells.Item(20, 1).Value = 46073.1875
$ws.Cells.Item(20, 2).Value = -34.352
$ws.Cells.Item(20, 3).Value = 633.444
$ws.Cells.Item(20, 4).Value = 633.444
$ws.Cells.Item(20, 6).Value = "20.02.202619"
$ws.Cells.Item(21, 1).Value = 46073.19791666666
$ws.Cells.Item(21, 2).Value = -95.896
$ws.Cells.Item(21, 3).Value = 1832.085
$ws.Cells.Item(21, 4).Value = 1832.085
$ws.Cells.Item(21, 6).Value = "20.02.202620"
$ws.Cells.Item(22, 1).Value = 46073.20833333334
$ws.Cells.Item(22, 2).Value = -43.445
$ws.Cells.Item(22, 3).Value = 2376.526
$ws.Cells.Item(22, 4).Value = 2376.526
$ws.Cells.Item(22, 6).Value = "20.02.202621"
$ws.Cells.Item(23, 1).Value = 46073.21875
$ws.Cells.Item(23, 2).Value = -96.169
$ws.Cells.Item(23, 3).Value = 2180.946
$ws.Cells.Item(23, 4).Value = 2180.946
$ws.Cells.Item(23, 6).Value = "20.02.202622"
$ws.Cells.Item(24, 1).Value = 46073.22916666666
$ws.Cells.Item(24, 2).Value = -114.308
$ws.Cells.Item(24, 3).Value = 1331.72
$ws.Cells.Item(24, 4).Value = 1331.72
$ws.Cells.Item(24, 6).Value = "20.02.202623"
$ws.Cells.Item(25, 1).Value = 46073.23958333334
$ws.Cells.Item(25, 2).Value = -120.568
$ws.Cells.Item(25, 3).Value = 4446.377
$ws.Cells.Item(25, 4).Value = 4446.377
$ws.Cells.Item(25, 6).Value = "20.02.202624"
$ws.Cells.Item(26, 1).Value = 46073.25
$ws.Cells.Item(26, 2).Value = -37.943
$ws.Cells.Item(26, 3).Value = 1240.118
$ws.Cells.Item(26, 4).Value = 1240.118
$ws.Cells.Item(26, 6).Value = "20.02.202625"
$ws.Cells.Item(27, 1).Value = 46073.26041666666
$ws.Cells.Item(27, 2).Value = -53.757
$ws.Cells.Item(27, 3).Value = 1176.126
$ws.Cells.Item(27, 4).Value = 1176.126
$ws.Cells.Item(27, 6).Value = "20.02.202626"
$ws.Cells.Item(28, 1).Value = 46073.27083333334
$ws.Cells.Item(28, 2).Value = -7.53
$ws.Cells.Item(28, 3).Value = 1202.142
$ws.Cells.Item(28, 4).Value = 1202.142
$ws.Cells.Item(28, 6).Value = "20.02.202627"
$ws.Cells.Item(29, 1).Value = 46073.28125
$ws.Cells.Item(29, 2).Value = -3.754
$ws.Cells.Item(29, 3).Value = 1205
$ws.Cells.Item(29, 4).Value = 1205
$ws.Cells.Item(29, 6).Value = "20.02.202628"
$ws.Cells.Item(30, 1).Value = 46073.29166666666
$ws.Cells.Item(30, 2).Value = -5.919
$ws.Cells.Item(30, 3).Value = 771.9640000000001
$ws.Cells.Item(30, 4).Value = 771.9640000000001
$ws.Cells.Item(30, 6).Value = "20.02.202629"
$ws.Cells.Item(31, 1).Value = 46073.30208333334
$ws.Cells.Item(31, 2).Value = 1.183
$ws.Cells.Item(31, 3).Value = 768.152
$ws.Cells.Item(31, 4).Value = 768.152
$ws.Cells.Item(31, 6).Value = "20.02.202630"
$ws.Cells.Item(32, 1).Value = 46073.3125
$ws.Cells.Item(32, 2).Value = 47.936
$ws.Cells.Item(32, 3).Value = 305.636
$ws.Cells.Item(32, 4).Value = 305.636
$ws.Cells.Item(32, 6).Value = "20.02.202631"
$ws.Cells.Item(33, 1).Value = 46073.32291666666
$ws.Cells.Item(33, 2).Value = 102.647
$ws.Cells.Item(33, 3).Value = 148.628
$ws.Cells.Item(33, 4).Value = 148.628
$ws.Cells.Item(33, 6).Value = "20.02.202632"
$ws.Cells.Item(34, 1).Value = 46073.33333333334
$ws.Cells.Item(34, 2).Value = 41.24
$ws.Cells.Item(34, 3).Value = 306.297
$ws.Cells.Item(34, 4).Value = 306.297
$ws.Cells.Item(34, 6).Value = "20.02.202633"
$ws.Cells.Item(35, 1).Value = 46073.34375
$ws.Cells.Item(35, 2).Value = 47.254
$ws.Cells.Item(35, 3).Value = 272.313
$ws.Cells.Item(35, 4).Value = 272.313
$ws.Cells.Item(35, 6).Value = "20.02.202634"
$ws.Cells.Item(36, 1).Value = 46073.35416666666
$ws.Cells.Item(36, 2).Value = 55.766
$ws.Cells.Item(36, 3).Value = -179.721
$ws.Cells.Item(36, 4).Value = -179.721
$ws.Cells.Item(36, 6).Value = "20.02.202635"
$ws.Cells.Item(37, 1).Value = 46073.36458333334
$ws.Cells.Item(37, 2).Value = 159.654
$ws.Cells.Item(37, 3).Value = -1228.219
$ws.Cells.Item(37, 4).Value = -1228.219
$ws.Cells.Item(37, 6).Value = "20.02.202636"
$ws.Cells.Item(38, 1).Value = 46073.375
$ws.Cells.Item(38, 2).Value = 65.68899999999999
$ws.Cells.Item(38, 3).Value = 18.253
$ws.Cells.Item(38, 4).Value = 18.253
$ws.Cells.Item(38, 6).Value = "20.02.202637"
$ws.Cells.Item(39, 1).Value = 46073.38541666666
$ws.Cells.Item(39, 2).Value = 93.417
$ws.Cells.Item(39, 3).Value = 269.834
$ws.Cells.Item(39, 4).Value = 269.834
$ws.Cells.Item(39, 6).Value = "20.02.202638"
$ws.Cells.Item(40, 1).Value = 46073.39583333334
$ws.Cells.Item(40, 2).Value = 77.452
$ws.Cells.Item(40, 3).Value = -117.5
$ws.Cells.Item(40, 4).Value = -117.5
$ws.Cells.Item(40, 6).Value = "20.02.202639"
$ws.Cells.Item(41, 1).Value = 46073.40625
$ws.Cells.Item(41, 2).Value = 84.142
$ws.Cells.Item(41, 3).Value = -32.037
$ws.Cells.Item(41, 4).Value = -32.037
$ws.Cells.Item(41, 6).Value = "20.02.202640"
$ws.Cells.Item(42, 1).Value = 46073.41666666666
$ws.Cells.Item(42, 2).Value = 17.198
$ws.Cells.Item(42, 3).Value = 18.134
$ws.Cells.Item(42, 4).Value = 18.134
$ws.Cells.Item(42, 6).Value = "20.02.202641"
$ws.Cells.Item(43, 1).Value = 46073.42708333334
$ws.Cells.Item(43, 2).Value = 2.272
$ws.Cells.Item(43, 3).Value = 4.731
$ws.Cells.Item(43, 4).Value = 4.731
$ws.Cells.Item(43, 6).Value = "20.02.202642"
$ws.Cells.Item(44, 1).Value = 46073.4375
$ws.Cells.Item(44, 2).Value = -40.457
$ws.Cells.Item(44, 3).Value = 601.886
$ws.Cells.Item(44, 4).Value = 601.886
$ws.Cells.Item(44, 6).Value = "20.02.202643"
$ws.Cells.Item(45, 1).Value = 46073.44791666666
$ws.Cells.Item(45, 2).Value = 18.399
$ws.Cells.Item(45, 3).Value = -649.668
$ws.Cells.Item(45, 4).Value = -649.668
$ws.Cells.Item(45, 6).Value = "20.02.202644"
$ws.Cells.Item(46, 1).Value = 46073.45833333334
$ws.Cells.Item(46, 2).Value = 29.361
$ws.Cells.Item(46, 3).Value = -821.3440000000001
$ws.Cells.Item(46, 4).Value = -821.3440000000001
$ws.Cells.Item(46, 6).Value = "20.02.202645"
$ws.Cells.Item(47, 1).Value = 46073.46875
$ws.Cells.Item(47, 2).Value = 42.625
$ws.Cells.Item(47, 3).Value = -1288.937
$ws.Cells.Item(47, 4).Value = -1288.937
$ws.Cells.Item(47, 6).Value = "20.02.202646"
$ws.Cells.Item(48, 1).Value = 46073.47916666666
$ws.Cells.Item(48, 2).Value = 63.291
$ws.Cells.Item(48, 3).Value = -3142.696
$ws.Cells.Item(48, 4).Value = -3142.696
$ws.Cells.Item(48, 6).Value = "20.02.202647"
$ws.Cells.Item(49, 1).Value = 46073.48958333334
$ws.Cells.Item(49, 2).Value = 91.333
$ws.Cells.Item(49, 3).Value = -3495.704
$ws.Cells.Item(49, 4).Value = -3495.704
$ws.Cells.Item(49, 6).Value = "20.02.202648"
$ws.Cells.Item(50, 1).Value = 46073.5
$ws.Cells.Item(50, 2).Value = 111.304
$ws.Cells.Item(50, 3).Value = -2752.498
$ws.Cells.Item(50, 4).Value = -2752.498
$ws.Cells.Item(50, 6).Value = "20.02.202649"
$ws.Cells.Item(51, 1).Value = 46073.51041666666
$ws.Cells.Item(51, 2).Value = 76.813
$ws.Cells.Item(51, 3).Value = -91.855
$ws.Cells.Item(51, 4).Value = -91.855
$ws.Cells.Item(51, 6).Value = "20.02.202650"
$ws.Cells.Item(52, 1).Value = 46073.52083333334
$ws.Cells.Item(52, 2).Value = 57.181
$ws.Cells.Item(52, 3).Value = -164.783
$ws.Cells.Item(52, 4).Value = -164.783
$ws.Cells.Item(52, 6).Value = "20.02.202651"
$ws.Cells.Item(53, 1).Value = 46073.53125
$ws.Cells.Item(53, 2).Value = 27.248
$ws.Cells.Item(53, 3).Value = -104.85
$ws.Cells.Item(53, 4).Value = -104.85
$ws.Cells.Item(53, 6).Value = "20.02.202652"
$ws.Cells.Item(54, 1).Value = 46073.54166666666
$ws.Cells.Item(54, 2).Value = 61.046
$ws.Cells.Item(54, 3).Value = 14.361
$ws.Cells.Item(54, 4).Value = 14.361
$ws.Cells.Item(54, 6).Value = "20.02.202653"
$ws.Cells.Item(55, 1).Value = 46073.55208333334
$ws.Cells.Item(55, 2).Value = 53.991
$ws.Cells.Item(55, 3).Value = 4.431
$ws.Cells.Item(55, 4).Value = 4.431
$ws.Cells.Item(55, 6).Value = "20.02.202654"
$ws.Cells.Item(56, 1).Value = 46073.5625
$ws.Cells.Item(56, 2).Value = 30.579
$ws.Cells.Item(56, 3).Value = 0.745
$ws.Cells.Item(56, 4).Value = 0.745
$ws.Cells.Item(56, 6).Value = "20.02.202655"
$ws.Cells.Item(57, 1).Value = 46073.57291666666
$ws.Cells.Item(57, 2).Value = 22.318
$ws.Cells.Item(57, 3).Value = 0.988
$ws.Cells.Item(57, 4).Value = 0.988
$ws.Cells.Item(57, 6).Value = "20.02.202656"
$ws.Cells.Item(58, 1).Value = 46073.58333333334
$ws.Cells.Item(58, 2).Value = 20.505
$ws.Cells.Item(58, 3).Value = 37.024
$ws.Cells.Item(58, 4).Value = 37.024
$ws.Cells.Item(58, 6).Value = "20.02.202657"
$ws.Cells.Item(59, 1).Value = 46073.59375
$ws.Cells.Item(59, 2).Value = 41.005
$ws.Cells.Item(59, 3).Value = 47.905
$ws.Cells.Item(59, 4).Value = 47.905
$ws.Cells.Item(59, 6).Value = "20.02.202658"
$ws.Cells.Item(60, 1).Value = 46073.60416666666
$ws.Cells.Item(60, 2).Value = 31.141
$ws.Cells.Item(60, 3).Value = -321.655
$ws.Cells.Item(60, 4).Value = -321.655
$ws.Cells.Item(60, 6).Value = "20.02.202659"
$ws.Cells.Item(61, 1).Value = 46073.61458333334
$ws.Cells.Item(61, 2).Value = 21.831
$ws.Cells.Item(61, 3).Value = 127.118
$ws.Cells.Item(61, 4).Value = 127.118
$ws.Cells.Item(61, 6).Value = "20.02.202660"
$ws.Cells.Item(62, 1).Value = 46073.625
$ws.Cells.Item(62, 2).Value = 64.142
$ws.Cells.Item(62, 3).Value = -2144.719
$ws.Cells.Item(62, 4).Value = -2144.719
$ws.Cells.Item(62, 6).Value = "20.02.202661"
$ws.Cells.Item(63, 1).Value = 46073.63541666666
$ws.Cells.Item(63, 2).Value = 18.557
$ws.Cells.Item(63, 3).Value = 4.906
$ws.Cells.Item(63, 4).Value = 4.906
$ws.Cells.Item(63, 6).Value = "20.02.202662"
$ws.Cells.Item(64, 1).Value = 46073.64583333334
$ws.Cells.Item(64, 2).Value = 59.884
$ws.Cells.Item(64, 3).Value = 299.842
$ws.Cells.Item(64, 4).Value = 299.842
$ws.Cells.Item(64, 6).Value = "20.02.202663"
$ws.Cells.Item(65, 1).Value = 46073.65625
$ws.Cells.Item(65, 2).Value = 56.032
$ws.Cells.Item(65, 3).Value = 59.822
$ws.Cells.Item(65, 4).Value = 59.822
$ws.Cells.Item(65, 6).Value = "20.02.202664"
$ws.Cells.Item(66, 1).Value = 46073.66666666666
$ws.Cells.Item(66, 2).Value = 72.998
$ws.Cells.Item(66, 3).Value = -497.501
$ws.Cells.Item(66, 4).Value = -497.501
$ws.Cells.Item(66, 6).Value = "20.02.202665"
$ws.Cells.Item(67, 1).Value = 46073.67708333334
$ws.Cells.Item(67, 2).Value = 30.151
$ws.Cells.Item(67, 3).Value = 13.995
$ws.Cells.Item(67, 4).Value = 13.995
$ws.Cells.Item(67, 6).Value = "20.02.202666"
$ws.Cells.Item(68, 1).Value = 46073.6875
$ws.Cells.Item(68, 2).Value = 15.786
$ws.Cells.Item(68, 3).Value = 0.3
$ws.Cells.Item(68, 4).Value = 0.3
$ws.Cells.Item(68, 6).Value = "20.02.202667"
$ws.Cells.Item(69, 1).Value = 46073.69791666666
$ws.Cells.Item(69, 2).Value = 8.286
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 6).Value = "20.02.202668"
$ws.Cells.Item(70, 1).Value = 46073.70833333334
$ws.Cells.Item(70, 2).Value = 50.022
$ws.Cells.Item(70, 3).Value = 272.254
$ws.Cells.Item(70, 4).Value = 272.254
$ws.Cells.Item(70, 6).Value = "20.02.202669"
$ws.Cells.Item(71, 1).Value = 46073.71875
$ws.Cells.Item(71, 2).Value = 54.116
$ws.Cells.Item(71, 3).Value = 165.444
$ws.Cells.Item(71, 4).Value = 165.444
$ws.Cells.Item(71, 6).Value = "20.02.202670"
$ws.Cells.Item(72, 1).Value = 46073.72916666666
$ws.Cells.Item(72, 2).Value = 61.226
$ws.Cells.Item(72, 3).Value = 269.43
$ws.Cells.Item(72, 4).Value = 269.43
$ws.Cells.Item(72, 6).Value = "20.02.202671"
$ws.Cells.Item(73, 1).Value = 46073.73958333334
$ws.Cells.Item(73, 2).Value = 77.601
$ws.Cells.Item(73, 3).Value = -64.39700000000001
$ws.Cells.Item(73, 4).Value = -64.39700000000001
$ws.Cells.Item(73, 6).Value = "20.02.202672"
$ws.Cells.Item(74, 1).Value = 46073.75
$ws.Cells.Item(74, 2).Value = 52.264
$ws.Cells.Item(74, 3).Value = -1418.26
$ws.Cells.Item(74, 4).Value = -1418.26
$ws.Cells.Item(74, 6).Value = "20.02.202673"
$ws.Cells.Item(75, 1).Value = 46073.76041666666
$ws.Cells.Item(75, 2).Value = 44.458
$ws.Cells.Item(75, 3).Value = -668.178
$ws.Cells.Item(75, 4).Value = -668.178
$ws.Cells.Item(75, 6).Value = "20.02.202674"
$ws.Cells.Item(76, 1).Value = 46073.77083333334
$ws.Cells.Item(76, 2).Value = 64.249
$ws.Cells.Item(76, 3).Value = -2157.106
$ws.Cells.Item(76, 4).Value = -2157.106
$ws.Cells.Item(76, 6).Value = "20.02.202675"
$ws.Cells.Item(77, 1).Value = 46073.78125
$ws.Cells.Item(77, 2).Value = 50.606
$ws.Cells.Item(77, 3).Value = -687.248
$ws.Cells.Item(77, 4).Value = -687.248
$ws.Cells.Item(77, 6).Value = "20.02.202676"
$ws.Cells.Item(78, 1).Value = 46073.79166666666
$ws.Cells.Item(78, 2).Value = 53.035
$ws.Cells.Item(78, 3).Value = -235.696
$ws.Cells.Item(78, 4).Value = -235.696
$ws.Cells.Item(78, 6).Value = "20.02.202677"
$ws.Cells.Item(79, 1).Value = 46073.80208333334
$ws.Cells.Item(79, 2).Value = 52.167
$ws.Cells.Item(79, 3).Value = -21.863
$ws.Cells.Item(79, 4).Value = -21.863
$ws.Cells.Item(79, 6).Value = "20.02.202678"
$ws.Cells.Item(80, 1).Value = 46073.8125
$ws.Cells.Item(80, 2).Value = 43.962
$ws.Cells.Item(80, 3).Value = 12.044
$ws.Cells.Item(80, 4).Value = 12.044
$ws.Cells.Item(80, 6).Value = "20.02.202679"
$ws.Cells.Item(81, 1).Value = 46073.82291666666
$ws.Cells.Item(81, 2).Value = 28.082
$ws.Cells.Item(81, 3).Value = 68.11499999999999
$ws.Cells.Item(81, 4).Value = 68.11499999999999
$ws.Cells.Item(81, 6).Value = "20.02.202680"
$ws.Cells.Item(82, 1).Value = 46073.83333333334
$ws.Cells.Item(82, 2).Value = 23.996
$ws.Cells.Item(82, 3).Value = -123.614
$ws.Cells.Item(82, 4).Value = -123.614
$ws.Cells.Item(82, 6).Value = "20.02.202681"
$ws.Cells.Item(83, 1).Value = 46073.84375
$ws.Cells.Item(83, 2).Value = 5.438
$ws.Cells.Item(83, 3).Value = 240.993
$ws.Cells.Item(83, 4).Value = 240.993
$ws.Cells.Item(83, 6).Value = "20.02.202682"
$ws.Cells.Item(84, 1).Value = 46073.85416666666
$ws.Cells.Item(84, 2).Value = 19.052
$ws.Cells.Item(84, 3).Value = 163.812
$ws.Cells.Item(84, 4).Value = 163.812
$ws.Cells.Item(84, 6).Value = "20.02.202683"
$ws.Cells.Item(85, 1).Value = 46073.86458333334
$ws.Cells.Item(85, 2).Value = 10.483
$ws.Cells.Item(85, 3).Value = 145.406
$ws.Cells.Item(85, 4).Value = 145.406
$ws.Cells.Item(85, 6).Value = "20.02.202684"
$ws.Cells.Item(86, 1).Value = 46073.875
$ws.Cells.Item(86, 2).Value = -12.545
$ws.Cells.Item(86, 3).Value = 599
$ws.Cells.Item(86, 4).Value = 599
$ws.Cells.Item(86, 6).Value = "20.02.202685"
$ws.Cells.Item(87, 1).Value = 46073.88541666666
$ws.Cells.Item(87, 2).Value = -54.078
$ws.Cells.Item(87, 3).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 6).Value = "20.02.202686"
$ws.Cells.Item(88, 1).Value = 46073.89583333334
$ws.Cells.Item(88, 2).Value = -48.358
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 6).Value = "20.02.202687"
$ws.Cells.Item(89, 1).Value = 46073.90625
$ws.Cells.Item(89, 2).Value = -36.799
$ws.Cells.Item(89, 3).Value = 628.235
$ws.Cells.Item(89, 4).Value = 628.235
$ws.Cells.Item(89, 6).Value = "20.02.202688"
$ws.Cells.Item(90, 1).Value = 46073.91666666666
$ws.Cells.Item(90, 2).Value = -4.304
$ws.Cells.Item(90, 3).Value = 228.916
$ws.Cells.Item(90, 4).Value = 228.916
$ws.Cells.Item(90, 6).Value = "20.02.202689"
$ws.Cells.Item(91, 1).Value = 46073.92708333334
$ws.Cells.Item(91, 2).Value = -17.421
$ws.Cells.Item(91, 3).Value = 613.2430000000001
$ws.Cells.Item(91, 4).Value = 613.2430000000001
$ws.Cells.Item(91, 6).Value = "20.02.202690"
$ws.Cells.Item(92, 1).Value = 46073.9375
$ws.Cells.Item(92, 2).Value = -29.436
$ws.Cells.Item(92, 3).Value = 620.133
$ws.Cells.Item(92, 4).Value = 620.133
$ws.Cells.Item(92, 6).Value = "20.02.202691"
$ws.Cells.Item(93, 1).Value = 46073.94791666666
$ws.Cells.Item(93, 2).Value = -2.224
$ws.Cells.Item(93, 3).Value = 599
$ws.Cells.Item(93, 4).Value = 599
$ws.Cells.Item(93, 6).Value = "20.02.202692"
$ws.Cells.Item(94, 1).Value = 46073.95833333334
$ws.Cells.Item(94, 2).Value = 11.951
$ws.Cells.Item(94, 3).Value = -101.021
$ws.Cells.Item(94, 4).Value = -101.021
$ws.Cells.Item(94, 6).Value = "20.02.202693"
$ws.Cells.Item(95, 1).Value = 46073.96875
$ws.Cells.Item(95, 2).Value = 24.191
$ws.Cells.Item(95, 3).Value = 59.409
$ws.Cells.Item(95, 4).Value = 59.409
$ws.Cells.Item(95, 6).Value = "20.02.202694"
$ws.Cells.Item(96, 1).Value = 46073.97916666666
$ws.Cells.Item(96, 2).Value = 13.736
$ws.Cells.Item(96, 3).Value = 19.675
$ws.Cells.Item(96, 4).Value = 19.675
$ws.Cells.Item(96, 6).Value = "20.02.202695"
$ws.Cells.Item(97, 1).Value = 46073.98958333334
$ws.Cells.Item(97, 2).Value = 14.049
$ws.Cells.Item(97, 3).Value = 106.554
$ws.Cells.Item(97, 4).Value = 106.554
$ws.Cells.Item(97, 6).Value = "20.02.202696"
$ws.Cells.Item(98, 1).Value = 46074
$ws.Cells.Item(98, 2).Value = 58.209
$ws.Cells.Item(98, 3).Value = -131.882
$ws.Cells.Item(98, 4).Value = -131.882
$ws.Cells.Item(98, 6).Value = "21.02.20261"
$ws.Cells.Item(99, 1).Value = 46074.01041666666
$ws.Cells.Item(99, 2).Value = 40.571
$ws.Cells.Item(99, 3).Value = 31.375
$ws.Cells.Item(99, 4).Value = 31.375
$ws.Cells.Item(99, 6).Value = "21.02.20262"
$ws.Cells.Item(100, 1).Value = 46074.02083333334
$ws.Cells.Item(100, 2).Value = -10.638
$ws.Cells.Item(100, 3).Value = 584.8680000000001
$ws.Cells.Item(100, 4).Value = 584.8680000000001
$ws.Cells.Item(100, 6).Value = "21.02.20263"
$ws.Cells.Item(101, 1).Value = 46074.03125
$ws.Cells.Item(101, 2).Value = -14.83
$ws.Cells.Item(101, 3).Value = 550
$ws.Cells.Item(101, 4).Value = 550
$ws.Cells.Item(101, 6).Value = "21.02.20264"
$ws.Cells.Item(102, 1).Value = 46074.04166666666
$ws.Cells.Item(102, 2).Value = -3.935
$ws.Cells.Item(102, 3).Value = 570
$ws.Cells.Item(102, 4).Value = 570
$ws.Cells.Item(102, 6).Value = "21.02.20265"
$ws.Cells.Item(103, 1).Value = 46074.05208333334
$ws.Cells.Item(103, 2).Value = -10.384
$ws.Cells.Item(103, 3).Value = 0
$ws.Cells.Item(103, 4).Value = 0
$ws.Cells.Item(103, 6).Value = "21.02.20266"
$ws.Cells.Item(104, 1).Value = 46074.0625
$ws.Cells.Item(104, 2).Value = -10.233
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = 0
$ws.Cells.Item(104, 6).Value = "21.02.20267"
$ws.Cells.Item(105, 1).Value = 46074.07291666666
$ws.Cells.Item(105, 2).Value = -7.787
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = 0
$ws.Cells.Item(105, 6).Value = "21.02.20268"
$ws.Cells.Item(106, 1).Value = 46074.08333333334
$ws.Cells.Item(106, 2).Value = 28.071
$ws.Cells.Item(106, 3).Value = 14.72
$ws.Cells.Item(106, 4).Value = 14.72
$ws.Cells.Item(106, 6).Value = "21.02.20269"
$ws.Cells.Item(107, 1).Value = 46074.09375
$ws.Cells.Item(107, 2).Value = 30.27
$ws.Cells.Item(107, 3).Value = -285.469
$ws.Cells.Item(107, 4).Value = -285.469
$ws.Cells.Item(107, 6).Value = "21.02.202610"
$ws.Cells.Item(108, 1).Value = 46074.10416666666
$ws.Cells.Item(108, 2).Value = 24.826
$ws.Cells.Item(108, 3).Value = 29.516
$ws.Cells.Item(108, 4).Value = 29.516
$ws.Cells.Item(108, 6).Value = "21.02.202611"
$ws.Cells.Item(109, 1).Value = 46074.11458333334
$ws.Cells.Item(109, 2).Value = 28.45
$ws.Cells.Item(109, 3).Value = -23.457
$ws.Cells.Item(109, 4).Value = -23.457
$ws.Cells.Item(109, 6).Value = "21.02.202612"
$ws.Cells.Item(110, 1).Value = 46074.125
$ws.Cells.Item(110, 2).Value = 30.718
$ws.Cells.Item(110, 3).Value = -7.306
$ws.Cells.Item(110, 4).Value = -7.306
$ws.Cells.Item(110, 6).Value = "21.02.202613"
$ws.Cells.Item(111, 1).Value = 46074.13541666666
$ws.Cells.Item(111, 2).Value = 39.956
$ws.Cells.Item(111, 3).Value = -347.661
$ws.Cells.Item(111, 4).Value = -347.661
$ws.Cells.Item(111, 6).Value = "21.02.202614"
$ws.Cells.Item(112, 1).Value = 46074.14583333334
$ws.Cells.Item(112, 2).Value = 23.402
$ws.Cells.Item(112, 3).Value = 11.216
$ws.Cells.Item(112, 4).Value = 11.216
$ws.Cells.Item(112, 6).Value = "21.02.202615"
$ws.Cells.Item(113, 1).Value = 46074.15625
$ws.Cells.Item(113, 2).Value = 6.848
$ws.Cells.Item(113, 3).Value = 90.89700000000001
$ws.Cells.Item(113, 4).Value = 90.89700000000001
$ws.Cells.Item(113, 6).Value = "21.02.202616"
$ws.Cells.Item(114, 1).Value = 46074.16666666666
$ws.Cells.Item(114, 2).Value = -19.556
$ws.Cells.Item(114, 3).Value = 601.7859999999999
$ws.Cells.Item(114, 4).Value = 601.7859999999999
$ws.Cells.Item(114, 6).Value = "21.02.202617"
$ws.Cells.Item(115, 1).Value = 46074.17708333334
$ws.Cells.Item(115, 2).Value = -45.263
$ws.Cells.Item(115, 3).Value = 622.149
$ws.Cells.Item(115, 4).Value = 622.149
$ws.Cells.Item(115, 6).Value = "21.02.202618"
$ws.Cells.Item(116, 1).Value = 46074.1875
$ws.Cells.Item(116, 2).Value = -82.938
$ws.Cells.Item(116, 3).Value = 683.905
$ws.Cells.Item(116, 4).Value = 683.905
$ws.Cells.Item(116, 6).Value = "21.02.202619"
$ws.Cells.Item(117, 1).Value = 46074.19791666666
$ws.Cells.Item(117, 2).Value = -133.656
$ws.Cells.Item(117, 3).Value = 602.816
$ws.Cells.Item(117, 4).Value = 602.816
$ws.Cells.Item(117, 6).Value = "21.02.202620"
$ws.Cells.Item(118, 1).Value = 46074.20833333334
$ws.Cells.Item(118, 2).Value = -96.083
$ws.Cells.Item(118, 3).Value = 2809.408
$ws.Cells.Item(118, 4).Value = 2809.408
$ws.Cells.Item(118, 6).Value = "21.02.202621"
$ws.Cells.Item(119, 1).Value = 46074.21875
$ws.Cells.Item(119, 2).Value = -92.59999999999999
$ws.Cells.Item(119, 3).Value = 1058.187
$ws.Cells.Item(119, 4).Value = 1058.187
$ws.Cells.Item(119, 6).Value = "21.02.202622"
$ws.Cells.Item(120, 1).Value = 46074.22916666666
$ws.Cells.Item(120, 2).Value = -66.952
$ws.Cells.Item(120, 3).Value = 817.6420000000001
$ws.Cells.Item(120, 4).Value = 817.6420000000001
$ws.Cells.Item(120, 6).Value = "21.02.202623"
$ws.Cells.Item(121, 1).Value = 46074.23958333334
$ws.Cells.Item(121, 2).Value = -66.55200000000001
$ws.Cells.Item(121, 3).Value = 815.298
$ws.Cells.Item(121, 4).Value = 815.298
$ws.Cells.Item(121, 6).Value = "21.02.202624"
$ws.Cells.Item(122, 1).Value = 46074.25
$ws.Cells.Item(122, 2).Value = -35.263
$ws.Cells.Item(122, 3).Value = 821.918
$ws.Cells.Item(122, 4).Value = 821.918
$ws.Cells.Item(122, 6).Value = "21.02.202625"
$ws.Cells.Item(123, 1).Value = 46074.26041666666
$ws.Cells.Item(123, 2).Value = -98.01000000000001
$ws.Cells.Item(123, 3).Value = 3691.703
$ws.Cells.Item(123, 4).Value = 3691.703
$ws.Cells.Item(123, 6).Value = "21.02.202626"
$ws.Cells.Item(124, 1).Value = 46074.27083333334
$ws.Cells.Item(124, 2).Value = -106.299
$ws.Cells.Item(124, 3).Value = 916.405
$ws.Cells.Item(124, 4).Value = 916.405
$ws.Cells.Item(124, 6).Value = "21.02.202627"
$ws.Cells.Item(125, 1).Value = 46074.28125
$ws.Cells.Item(125, 2).Value = -173.24
$ws.Cells.Item(125, 3).Value = 817.173
$ws.Cells.Item(125, 4).Value = 817.173
$ws.Cells.Item(125, 6).Value = "21.02.202628"
$ws.Cells.Item(126, 1).Value = 46074.29166666666
$ws.Cells.Item(126, 2).Value = -168.342
$ws.Cells.Item(126, 3).Value = 862.967
$ws.Cells.Item(126, 4).Value = 862.967
$ws.Cells.Item(126, 6).Value = "21.02.202629"
$ws.Cells.Item(127, 1).Value = 46074.30208333334
$ws.Cells.Item(127, 2).Value = 0
$ws.Cells.Item(127, 3).Value = 0
$ws.Cells.Item(127, 4).Value = 0
$ws.Cells.Item(127, 6).Value = "21.02.202630"
$ws.Cells.Item(128, 1).Value = 46074.3125
$ws.Cells.Item(128, 2).Value = 0
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(128, 4).Value = 0
$ws.Cells.Item(128, 6).Value = "21.02.202631"
$ws.Cells.Item(129, 1).Value = 46074.32291666666
$ws.Cells.Item(129, 2).Value = 0
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 0
$ws.Cells.Item(129, 6).Value = "21.02.202632"
$ws.Cells.Item(130, 1).Value = 46074.33333333334
$ws.Cells.Item(130, 2).Value = 0
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 0
$ws.Cells.Item(130, 6).Value = "21.02.202633"
$ws.Cells.Item(131, 1).Value = 46074.34375
$ws.Cells.Item(131, 2).Value = 0
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 0
$ws.Cells.Item(131, 6).Value = "21.02.202634"
$ws.Cells.Item(132, 1).Value = 46074.35416666666
$ws.Cells.Item(132, 2).Value = 0
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = 0
$ws.Cells.Item(132, 6).Value = "21.02.202635"
$ws.Cells.Item(133, 1).Value = 46074.36458333334
$ws.Cells.Item(133, 2).Value = 0
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(133, 4).Value = 0
$ws.Cells.Item(133, 6).Value = "21.02.202636"
$ws.Cells.Item(134, 1).Value = 46074.375
$ws.Cells.Item(134, 2).Value = 0
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(134, 4).Value = 0
$ws.Cells.Item(134, 6).Value = "21.02.202637"
$ws.Cells.Item(135, 1).Value = 46074.38541666666
$ws.Cells.Item(135, 2).Value = 0
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 0
$ws.Cells.Item(135, 6).Value = "21.02.202638"
$ws.Cells.Item(136, 1).Value = 46074.39583333334
$ws.Cells.Item(136, 2).Value = 0
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = 0
$ws.Cells.Item(136, 6).Value = "21.02.202639"
$ws.Cells.Item(137, 1).Value = 46074.40625
$ws.Cells.Item(137, 2).Value = 0
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = 0
$ws.Cells.Item(137, 6).Value = "21.02.202640"
$ws.Cells.Item(138, 1).Value = 46074.41666666666
$ws.Cells.Item(138, 2).Value = 0
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 0
$ws.Cells.Item(138, 6).Value = "21.02.202641"
$ws.Cells.Item(139, 1).Value = 46074.42708333334
$ws.Cells.Item(139, 2).Value = 0
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(139, 4).Value = 0
$ws.Cells.Item(139, 6).Value = "21.02.202642"
$ws.Cells.Item(140, 1).Value = 46074.4375
$ws.Cells.Item(140, 2).Value = 0
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(140, 4).Value = 0
$ws.Cells.Item(140, 6).Value = "21.02.202643"
$ws.Cells.Item(141, 1).Value = 46074.44791666666
$ws.Cells.Item(141, 2).Value = 0
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = 0
$ws.Cells.Item(141, 6).Value = "21.02.202644"
$ws.Cells.Item(142, 1).Value = 46074.45833333334
$ws.Cells.Item(142, 2).Value = 0
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(142, 4).Value = 0
$ws.Cells.Item(142, 6).Value = "21.02.202645"
$ws.Cells.Item(143, 1).Value = 46074.46875
$ws.Cells.Item(143, 2).Value = 0
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 6).Value = "21.02.202646"
$ws.Cells.Item(144, 1).Value = 46074.47916666666
$ws.Cells.Item(144, 2).Value = 0
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 0
$ws.Cells.Item(144, 6).Value = "21.02.202647"
$ws.Cells.Item(145, 1).Value = 46074.48958333334
$ws.Cells.Item(145, 2).Value = 0
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 0
$ws.Cells.Item(145, 6).Value = "21.02.202648"
$ws.Cells.Item(146, 1).Value = 46074.5
$ws.Cells.Item(146, 2).Value = 0
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 0
$ws.Cells.Item(146, 6).Value = "21.02.202649"
$ws.Cells.Item(147, 1).Value = 46074.51041666666
$ws.Cells.Item(147, 2).Value = 0
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 0
$ws.Cells.Item(147, 6).Value = "21.02.202650"
$ws.Cells.Item(148, 1).Value = 46074.52083333334
$ws.Cells.Item(148, 2).Value = 0
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 0
$ws.Cells.Item(148, 6).Value = "21.02.202651"
$ws.Cells.Item(149, 1).Value = 46074.53125
$ws.Cells.Item(149, 2).Value = 0
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 0
$ws.Cells.Item(149, 6).Value = "21.02.202652"
$ws.Cells.Item(150, 1).Value = 46074.54166666666
$ws.Cells.Item(150, 2).Value = 0
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 0
$ws.Cells.Item(150, 6).Value = "21.02.202653"
$ws.Cells.Item(151, 1).Value = 46074.55208333334
$ws.Cells.Item(151, 2).Value = 0
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 0
$ws.Cells.Item(151, 6).Value = "21.02.202654"
$ws.Cells.Item(152, 1).Value = 46074.5625
$ws.Cells.Item(152, 2).Value = 0
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(152, 4).Value = 0
$ws.Cells.Item(152, 6).Value = "21.02.202655"
$ws.Cells.Item(153, 1).Value = 46074.57291666666
$ws.Cells.Item(153, 2).Value = 0
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(153, 4).Value = 0
$ws.Cells.Item(153, 6).Value = "21.02.202656"
$ws.Cells.Item(154, 1).Value = 46074.58333333334
$ws.Cells.Item(154, 2).Value = 0
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(154, 4).Value = 0
$ws.Cells.Item(154, 6).Value = "21.02.202657"
$ws.Cells.Item(155, 1).Value = 46074.59375
$ws.Cells.Item(155, 2).Value = 0
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 0
$ws.Cells.Item(155, 6).Value = "21.02.202658"
$ws.Cells.Item(156, 1).Value = 46074.60416666666
$ws.Cells.Item(156, 2).Value = 0
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 0
$ws.Cells.Item(156, 6).Value = "21.02.202659"
$ws.Cells.Item(157, 1).Value = 46074.61458333334
$ws.Cells.Item(157, 2).Value = 0
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 0
$ws.Cells.Item(157, 6).Value = "21.02.202660"
$ws.Cells.Item(158, 1).Value = 46074.625
$ws.Cells.Item(158, 2).Value = 0
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 6).Value = "21.02.202661"
$ws.Cells.Item(159, 1).Value = 46074.63541666666
$ws.Cells.Item(159, 2).Value = 0
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 6).Value = "21.02.202662"
$ws.Cells.Item(160, 1).Value = 46074.64583333334
$ws.Cells.Item(160, 2).Value = 0
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 0
$ws.Cells.Item(160, 6).Value = "21.02.202663"
$ws.Cells.Item(161, 1).Value = 46074.65625
$ws.Cells.Item(161, 2).Value = 0
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 0
$ws.Cells.Item(161, 6).Value = "21.02.202664"
$ws.Cells.Item(162, 1).Value = 46074.66666666666
$ws.Cells.Item(162, 2).Value = 0
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 0
$ws.Cells.Item(162, 6).Value = "21.02.202665"
$ws.Cells.Item(163, 1).Value = 46074.67708333334
$ws.Cells.Item(163, 2).Value = 0
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(163, 4).Value = 0
$ws.Cells.Item(163, 6).Value = "21.02.202666"
$ws.Cells.Item(164, 1).Value = 46074.6875
$ws.Cells.Item(164, 2).Value = 0
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 6).Value = "21.02.202667"
$ws.Cells.Item(165, 1).Value = 46074.69791666666
$ws.Cells.Item(165, 2).Value = 0
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 6).Value = "21.02.202668"
$ws.Cells.Item(166, 1).Value = 46074.70833333334
$ws.Cells.Item(166, 2).Value = 0
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 0
$ws.Cells.Item(166, 6).Value = "21.02.202669"
$ws.Cells.Item(167, 1).Value = 46074.71875
$ws.Cells.Item(167, 2).Value = 0
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 0
$ws.Cells.Item(167, 6).Value = "21.02.202670"
$ws.Cells.Item(168, 1).Value = 46074.72916666666
$ws.Cells.Item(168, 2).Value = 0
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 0
$ws.Cells.Item(168, 6).Value = "21.02.202671"
$ws.Cells.Item(169, 1).Value = 46074.73958333334
$ws.Cells.Item(169, 2).Value = 0
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(169, 4).Value = 0
$ws.Cells.Item(169, 6).Value = "21.02.202672"
$ws.Cells.Item(170, 1).Value = 46074.75
$ws.Cells.Item(170, 2).Value = 0
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 6).Value = "21.02.202673"
$ws.Cells.Item(171, 1).Value = 46074.76041666666
$ws.Cells.Item(171, 2).Value = 0
$ws.Cells.Item(171, 3).Value = 0
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 6).Value = "21.02.202674"
$ws.Cells.Item(172, 1).Value = 46074.77083333334
$ws.Cells.Item(172, 2).Value = 0
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 6).Value = "21.02.202675"
$ws.Cells.Item(173, 1).Value = 46074.78125
$ws.Cells.Item(173, 2).Value = 0
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 6).Value = "21.02.202676"
$ws.Cells.Item(174, 1).Value = 46074.79166666666
$ws.Cells.Item(174, 2).Value = 0
$ws.Cells.Item(174, 3).Value = 0
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 6).Value = "21.02.202677"
$ws.Cells.Item(175, 1).Value = 46074.80208333334
$ws.Cells.Item(175, 2).Value = 0
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 6).Value = "21.02.202678"
$ws.Cells.Item(176, 1).Value = 46074.8125
$ws.Cells.Item(176, 2).Value = 0
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 0
$ws.Cells.Item(176, 6).Value = "21.02.202679"
$ws.Cells.Item(177, 1).Value = 46074.82291666666
$ws.Cells.Item(177, 2).Value = 0
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 0
$ws.Cells.Item(177, 6).Value = "21.02.202680"
$ws.Cells.Item(178, 1).Value = 46074.83333333334
$ws.Cells.Item(178, 2).Value = 0
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 0
$ws.Cells.Item(178, 6).Value = "21.02.202681"
$ws.Cells.Item(179, 1).Value = 46074.84375
$ws.Cells.Item(179, 2).Value = 0
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 0
$ws.Cells.Item(179, 6).Value = "21.02.202682"
$ws.Cells.Item(180, 1).Value = 46074.85416666666
$ws.Cells.Item(180, 2).Value = 0
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 6).Value = "21.02.202683"
$ws.Cells.Item(181, 1).Value = 46074.86458333334
$ws.Cells.Item(181, 2).Value = 0
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 6).Value = "21.02.202684"
$ws.Cells.Item(182, 1).Value = 46074.875
$ws.Cells.Item(182, 2).Value = 0
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 6).Value = "21.02.202685"
$ws.Cells.Item(183, 1).Value = 46074.88541666666
$ws.Cells.Item(183, 2).Value = 0
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 6).Value = "21.02.202686"
$ws.Cells.Item(184, 1).Value = 46074.89583333334
$ws.Cells.Item(184, 2).Value = 0
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 6).Value = "21.02.202687"
$ws.Cells.Item(185, 1).Value = 46074.90625
$ws.Cells.Item(185, 2).Value = 0
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 0
$ws.Cells.Item(185, 6).Value = "21.02.202688"
$ws.Cells.Item(186, 1).Value = 46074.91666666666
$ws.Cells.Item(186, 2).Value = 0
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 6).Value = "21.02.202689"
$ws.Cells.Item(187, 1).Value = 46074.92708333334
$ws.Cells.Item(187, 2).Value = 0
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = 0
$ws.Cells.Item(187, 6).Value = "21.02.202690"
$ws.Cells.Item(188, 1).Value = 46074.9375
$ws.Cells.Item(188, 2).Value = 0
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 6).Value = "21.02.202691"
$ws.Cells.Item(189, 1).Value = 46074.94791666666
$ws.Cells.Item(189, 2).Value = 0
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 0
$ws.Cells.Item(189, 6).Value = "21.02.202692"
$ws.Cells.Item(190, 1).Value = 46074.95833333334
$ws.Cells.Item(190, 2).Value = 0
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(190, 4).Value = 0
$ws.Cells.Item(190, 6).Value = "21.02.202693"
$ws.Cells.Item(191, 1).Value = 46074.96875
$ws.Cells.Item(191, 2).Value = 0
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 6).Value = "21.02.202694"
$ws.Cells.Item(192, 1).Value = 46074.97916666666
$ws.Cells.Item(192, 2).Value = 0
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 6).Value = "21.02.202695"
$ws.Cells.Item(193, 1).Value = 46074.98958333334
$ws.Cells.Item(193, 2).Value = 0
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 0
$ws.Cells.Item(193, 6).Value = "21.02.202696"
